$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 105
$ws.Range("F3").Value = 74
$ws.Range("F5").Value = 1274
$ws.Range("F6").Value = 1764
$ws.Range("F8").Value = 582
$ws.Range("F9").Value = 2712
$ws.Range("F10").Value = 754
$ws.Range("F11").Value = 578
$ws.Range("F13").Value = 54
$ws.Range("F15").Value = 360
$ws.Range("F16").Value = 360
$ws.Range("F20").Value = 719
$ws.Range("F21").Value = 11
$ws.Range("F28").Value = 497
$ws.Range("F31").Value = 536
$ws.Range("F32").Value = 555
$ws.Range("F35").Value = 355
$ws.Range("F36").Value = 4615
$ws.Range("F37").Value = 183

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4217
$ws.Range("F10").Value = 321
$ws.Range("F11").Value = 338
$ws.Range("F15").Value = 159
$ws.Range("F17").Value = 275
$ws.Range("F22").Value = 1766
$ws.Range("F24").Value = 265
$ws.Range("F25").Value = 17
$ws.Range("F28").Value = 7
$ws.Range("F34").Value = 8
$ws.Range("F35").Value = 488

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 552
$ws.Range("F7").Value = 166
$ws.Range("F8").Value = 218

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 552
$ws.Range("F4").Value = 105
$ws.Range("F5").Value = 74
$ws.Range("F7").Value = 1274
$ws.Range("F8").Value = 1764
$ws.Range("F10").Value = 166
$ws.Range("F12").Value = 582
$ws.Range("F13").Value = 2712
$ws.Range("F14").Value = 754
$ws.Range("F17").Value = 360
$ws.Range("F18").Value = 360
$ws.Range("F19").Value = 321
$ws.Range("F20").Value = 338
$ws.Range("F24").Value = 719
$ws.Range("F26").Value = 159
$ws.Range("F27").Value = 11
$ws.Range("F31").Value = 275
$ws.Range("F34").Value = 218
$ws.Range("F37").Value = 1766
$ws.Range("F38").Value = 497
$ws.Range("F39").Value = 265
$ws.Range("F40").Value = 17
$ws.Range("F44").Value = 355
$ws.Range("F45").Value = 4615
$ws.Range("F46").Value = 183
$ws.Range("F48").Value = 488
